$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 42.675953
$ws.Range("H2").Value = 128.027859
$ws.Range("I2").Value = 0.9123907898314253
$ws.Range("J2").Value = 0.9123907898314252
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02725333333333333
$ws.Range("N2").Value = 0.08176
$ws.Range("O2").Value = 0.0007089206372884383
$ws.Range("P2").Value = 0.0007089206372884382
$ws.Range("Q2").Value = 1.163061972426667
$ws.Range("R2").Value = 10.46755775184
$ws.Range("S2").Value = 0.0006468126601833956
$ws.Range("T2").Value = 0.0006468126601833954
$ws.Range("G3").Value = 42.675953
$ws.Range("H3").Value = 128.027859
$ws.Range("I3").Value = 0.9123907898314253
$ws.Range("J3").Value = 0.9123907898314252
$ws.Range("M3").Value = 38.416166
$ws.Range("N3").Value = 115.248498
$ws.Range("O3").Value = 0.9992910793627116
$ws.Range("P3").Value = 0.9992910793627116
$ws.Range("Q3").Value = 1639.446494656198
$ws.Range("R3").Value = 14755.01845190578
$ws.Range("S3").Value = 0.9117439771712419
$ws.Range("T3").Value = 0.9117439771712418
$ws.Range("I4").Value = 0.03270518515803058
$ws.Range("J4").Value = 0.03270518515803057
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02725333333333333
$ws.Range("N4").Value = 0.08176
$ws.Range("O4").Value = 0.0007089206372884383
$ws.Range("P4").Value = 0.0007089206372884382
$ws.Range("Q4").Value = 0.04169064131555556
$ws.Range("R4").Value = 0.37521577184
$ws.Range("S4").Value = 0.00002318538070486741
$ws.Range("T4").Value = 0.0000231853807048674
$ws.Range("I5").Value = 0.03270518515803058
$ws.Range("J5").Value = 0.03270518515803057
$ws.Range("M5").Value = 38.416166
$ws.Range("N5").Value = 115.248498
$ws.Range("O5").Value = 0.9992910793627116
$ws.Range("P5").Value = 0.9992910793627116
$ws.Range("Q5").Value = 58.76692505228133
$ws.Range("R5").Value = 528.902325470532
$ws.Range("S5").Value = 0.03268199977732571
$ws.Range("T5").Value = 0.03268199977732571
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9319226666666666
$ws.Range("H6").Value = 2.795768
$ws.Range("I6").Value = 0.01992404616955614
$ws.Range("J6").Value = 0.01992404616955614
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02725333333333333
$ws.Range("N6").Value = 0.08176
$ws.Range("O6").Value = 0.0007089206372884383
$ws.Range("P6").Value = 0.0007089206372884382
$ws.Range("Q6").Value = 0.02539799907555555
$ws.Range("R6").Value = 0.22858199168
$ws.Range("S6").Value = 0.00001412456750788601
$ws.Range("T6").Value = 0.00001412456750788601
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9319226666666666
$ws.Range("H7").Value = 2.795768
$ws.Range("I7").Value = 0.01992404616955614
$ws.Range("J7").Value = 0.01992404616955614
$ws.Range("M7").Value = 38.416166
$ws.Range("N7").Value = 115.248498
$ws.Range("O7").Value = 0.9992910793627116
$ws.Range("P7").Value = 0.9992910793627116
$ws.Range("Q7").Value = 35.80089586182932
$ws.Range("R7").Value = 322.208062756464
$ws.Range("S7").Value = 0.01990992160204826
$ws.Range("T7").Value = 0.01990992160204826
$ws.Range("G8").Value = 1.636145333333333
$ws.Range("H8").Value = 4.908436
$ws.Range("I8").Value = 0.03497997884098805
$ws.Range("J8").Value = 0.03497997884098804
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02725333333333333
$ws.Range("N8").Value = 0.08176
$ws.Range("O8").Value = 0.0007089206372884383
$ws.Range("P8").Value = 0.0007089206372884382
$ws.Range("Q8").Value = 0.04459041415111111
$ws.Range("R8").Value = 0.40131372736
$ws.Range("S8").Value = 0.00002479802889228934
$ws.Range("T8").Value = 0.00002479802889228933
$ws.Range("G9").Value = 1.636145333333333
$ws.Range("H9").Value = 4.908436
$ws.Range("I9").Value = 0.03497997884098805
$ws.Range("J9").Value = 0.03497997884098804
$ws.Range("M9").Value = 38.416166
$ws.Range("N9").Value = 115.248498
$ws.Range("O9").Value = 0.9992910793627116
$ws.Range("P9").Value = 0.9992910793627116
$ws.Range("Q9").Value = 62.85443072545866
$ws.Range("R9").Value = 565.689876529128
$ws.Range("S9").Value = 0.03495518081209576
$ws.Range("T9").Value = 0.03495518081209575
